# ADNIsubjects_fmri_PreprocessTrack.xlsx - "Add files via upload"
#
# The underlying edit fills in missing "Y" (and one "X") tracking marks in
# columns B/C/D of Sheet1 for subjects whose processing status had not yet
# been recorded, plus it moves the user's active selection/scroll position.
# The J16/J18/J20 COUNTIF/COUNTA summary cells recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 48: column C already says "Y"; column D was still blank -> "Y"
$ws.Range("D48").Value = "Y"

# Rows 85-89: columns C and D were blank -> "Y" (same as existing column B)
foreach ($r in 85..89) {
    $ws.Cells.Item($r, 3).Value = "Y"
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Row 90: only column D was blank, and it is marked "X" (not "Y")
$ws.Range("D90").Value = "X"

# Rows 91-94: columns C and D were blank -> "Y"
foreach ($r in 91..94) {
    $ws.Cells.Item($r, 3).Value = "Y"
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Rows 98-110: column B was blank -> "Y"
foreach ($r in 98..110) {
    $ws.Cells.Item($r, 2).Value = "Y"
}

# Update the saved view state: scroll position and active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 82
$win.ScrollColumn = 1
$ws.Range("B111").Select()
